$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 51980.668
$ws.Range("I21").Value = 11750.5
$ws.Range("J21").Value = 63475
$ws.Range("K21").Value = 11750.5
$ws.Range("L21").Value = 63475
$ws.Range("M21").Value = -11282.5
$ws.Range("N21").Value = -64411
$ws.Range("H23").Value = 51980.668
$ws.Range("I23").Value = 11750.5
$ws.Range("J23").Value = 63475
$ws.Range("K23").Value = 11750.5
$ws.Range("L23").Value = 63475
$ws.Range("M23").Value = -11516.5
$ws.Range("N23").Value = -63943
$ws.Range("H38").Value = 34.363636
$ws.Range("I38").Value = 34.363636
$ws.Range("K38").Value = 103.090908
$ws.Range("M38").Value = 268.909092
$ws.Range("H39").Value = 160.53334
$ws.Range("I39").Value = 85.75
$ws.Range("K39").Value = 257.25
$ws.Range("M39").Value = 38.75
$ws.Range("H40").Value = 4389.25
$ws.Range("I40").Value = 3929.4167
$ws.Range("J40").Value = 4586.3213
$ws.Range("K40").Value = 3929.4167
$ws.Range("L40").Value = 4586.3213
$ws.Range("M40").Value = -3754.4167
$ws.Range("N40").Value = -4936.3213
$ws.Range("H70").Value = 7079.476
$ws.Range("J70").Value = 7130.1953
$ws.Range("L70").Value = 21390.5859
$ws.Range("N70").Value = -21930.5859
$ws.Range("H73").Value = 7079.476
$ws.Range("J73").Value = 7130.1953
$ws.Range("L73").Value = 21390.5859
$ws.Range("N73").Value = -23262.5859
$ws.Range("H76").Value = 7763.3335
$ws.Range("I76").Value = 7650
$ws.Range("K76").Value = 7650
$ws.Range("M76").Value = -7335
$ws.Range("H79").Value = 7763.3335
$ws.Range("I79").Value = 7650
$ws.Range("K79").Value = 7650
$ws.Range("M79").Value = -6558
$ws.Range("H96").Value = 475.25
$ws.Range("I96").Value = 400.2857
$ws.Range("K96").Value = 1200.8571
$ws.Range("M96").Value = 172.1428999999998
$ws.Range("H98").Value = 1010.0645
$ws.Range("I98").Value = 843.73334
$ws.Range("K98").Value = 843.73334
$ws.Range("M98").Value = 654.26666
$ws.Range("H122").Value = 1010.0645
$ws.Range("I122").Value = 843.73334
$ws.Range("K122").Value = 2531.20002
$ws.Range("M122").Value = -81.20002000000022
$ws.Range("H137").Value = 2623.8445
$ws.Range("I137").Value = 2239.04
$ws.Range("J137").Value = 3104.85
$ws.Range("K137").Value = 6717.12
$ws.Range("L137").Value = 9314.549999999999
$ws.Range("M137").Value = -4167.12
$ws.Range("N137").Value = -14414.55

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 6942.4375
$ws.Range("H45").Value = 51754.047
$ws.Range("I45").Value = 85288.914
$ws.Range("K45").Value = 85288.914
$ws.Range("M45").Value = -84911.914
$ws.Range("H55").Value = 50262
$ws.Range("J55").Value = 53000
$ws.Range("L55").Value = 53000
$ws.Range("N55").Value = -53630
$ws.Range("H86").Value = 8284
$ws.Range("I86").Value = 8284
$ws.Range("K86").Value = 8284
$ws.Range("M86").Value = -7098
$ws.Range("H89").Value = 8284
$ws.Range("I89").Value = 8284
$ws.Range("K89").Value = 24852
$ws.Range("M89").Value = -18924
$ws.Range("H132").Value = 2101.0908
$ws.Range("I132").Value = 1662.04
$ws.Range("K132").Value = 4986.12
$ws.Range("M132").Value = -2456.12

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4772
$ws.Range("J94").Value = 14117.385
$ws.Range("L94").Value = 14117.385
$ws.Range("N94").Value = -15019.385
$ws.Range("H134").Value = 2563.2222
$ws.Range("I134").Value = 780.14703
$ws.Range("J134").Value = 8074.5454
$ws.Range("K134").Value = 2340.44109
$ws.Range("L134").Value = 24223.6362
$ws.Range("M134").Value = 194.5589100000002
$ws.Range("N134").Value = -29293.6362

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3649.5588
$ws.Range("I31").Value = 2470.7
$ws.Range("K31").Value = 2470.7
$ws.Range("M31").Value = -2175.7
$ws.Range("H34").Value = 3649.5588
$ws.Range("I34").Value = 2470.7
$ws.Range("K34").Value = 2470.7
$ws.Range("M34").Value = -2268.7
$ws.Range("H132").Value = 52359.4
$ws.Range("I132").Value = 54862.58
$ws.Range("K132").Value = 164587.74
$ws.Range("M132").Value = -162057.74
$ws.Range("H134").Value = 22768.8
$ws.Range("I134").Value = 30237.281
$ws.Range("K134").Value = 90711.84299999999
$ws.Range("M134").Value = -88176.84299999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 120.77778
$ws.Range("J12").Value = 99.75
$ws.Range("L12").Value = 299.25
$ws.Range("N12").Value = -645.25
$ws.Range("H56").Value = 10422993
$ws.Range("I56").Value = 10422993
$ws.Range("K56").Value = 10422993
$ws.Range("M56").Value = -10422463
$ws.Range("H137").Value = 3565.2354
$ws.Range("J137").Value = 4955.1
$ws.Range("L137").Value = 14865.3
$ws.Range("N137").Value = -25065.3
$ws.Range("H138").Value = 3238.125
$ws.Range("I138").Value = 2843.5715
$ws.Range("K138").Value = 8530.7145
$ws.Range("M138").Value = -3390.7145

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 111921.766
$ws.Range("I122").Value = 187467.1
$ws.Range("K122").Value = 562401.3
$ws.Range("M122").Value = -559951.3
$ws.Range("H126").Value = 3630.6428
$ws.Range("I126").Value = 3404.0833
$ws.Range("K126").Value = 10212.2499
$ws.Range("M126").Value = -7742.249899999999
$ws.Range("H132").Value = 4959
$ws.Range("I132").Value = 4537.9414
$ws.Range("K132").Value = 13613.8242
$ws.Range("M132").Value = -11083.8242

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8616.666999999999
$ws.Range("I7").Value = 5610.8
$ws.Range("J7").Value = 12374
$ws.Range("K7").Value = 5610.8
$ws.Range("L7").Value = 12374
$ws.Range("M7").Value = -5498.8
$ws.Range("N7").Value = -12598
$ws.Range("H16").Value = 1243.3572
$ws.Range("I16").Value = 1072.238
$ws.Range("J16").Value = 1756.7142
$ws.Range("K16").Value = 1072.238
$ws.Range("L16").Value = 1756.7142
$ws.Range("M16").Value = -902.2380000000001
$ws.Range("N16").Value = -2096.7142
$ws.Range("H40").Value = 10343.25
$ws.Range("J40").Value = 12949.2
$ws.Range("L40").Value = 12949.2
$ws.Range("N40").Value = -13221.2
$ws.Range("H46").Value = 6688.2593
$ws.Range("J46").Value = 7658.6924
$ws.Range("L46").Value = 7658.6924
$ws.Range("N46").Value = -8034.6924
$ws.Range("H93").Value = 2675.76
$ws.Range("I93").Value = 2416.9167
$ws.Range("K93").Value = 2416.9167
$ws.Range("M93").Value = -1168.9167
$ws.Range("H122").Value = 5088.387
$ws.Range("J122").Value = 7351.364
$ws.Range("L122").Value = 22054.092
$ws.Range("N122").Value = -26954.092
$ws.Range("H126").Value = 8616.666999999999
$ws.Range("I126").Value = 5610.8
$ws.Range("J126").Value = 12374
$ws.Range("K126").Value = 16832.4
$ws.Range("L126").Value = 37122
$ws.Range("M126").Value = -14362.4
$ws.Range("N126").Value = -42062
$ws.Range("H132").Value = 6613.829
$ws.Range("I132").Value = 6764.5
$ws.Range("K132").Value = 20293.5
$ws.Range("M132").Value = -17763.5
$ws.Range("H136").Value = 47115.688
$ws.Range("I136").Value = 55957.73
$ws.Range("J136").Value = 6221.25
$ws.Range("K136").Value = 167873.19
$ws.Range("L136").Value = 18663.75
$ws.Range("M136").Value = -165323.19
$ws.Range("N136").Value = -23763.75
$ws.Range("H140").Value = 100214
$ws.Range("J140").Value = 100214
$ws.Range("L140").Value = 100214
$ws.Range("N140").Value = -110574

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H41").Value = 77531.14
$ws.Range("J41").Value = 77531.14
$ws.Range("L41").Value = 77531.14
$ws.Range("N41").Value = -78311.14
$ws.Range("H132").Value = 464185.1
$ws.Range("I132").Value = 12875.154
$ws.Range("J132").Value = 1116077.2
$ws.Range("K132").Value = 38625.462
$ws.Range("L132").Value = 3348231.6
$ws.Range("M132").Value = -36095.462
$ws.Range("N132").Value = -3353291.6
